$d = $word.ActiveDocument

# Helper: replace the visible text of a table cell, preserving the
# trailing end-of-cell marker (the cell Range includes it as one unit).
function Set-CellText($table, $row, $col, $newText) {
    $cellRange = $table.Cell($row, $col).Range
    $textRange = $d.Range($cellRange.Start, $cellRange.End - 1)
    $textRange.Text = $newText
}

# Update the date heading (first paragraph), preserving the paragraph mark.
$dateRange = $d.Paragraphs.Item(1).Range
$dateTextRange = $d.Range($dateRange.Start, $dateRange.End - 1)
$dateTextRange.Text = "2023-09-24 Sunday"

$t = $d.Tables.Item(1)

# Row 1 (0-based data row 0)
Set-CellText $t 1 1 "68÷6="   # was "84÷6="
Set-CellText $t 1 2 "95÷4="   # was "83÷9="
Set-CellText $t 1 3 "84÷7="   # was "61÷2="
Set-CellText $t 1 4 "21÷3="   # was "14÷4="
Set-CellText $t 1 5 "17÷5="   # was "25÷2="

# Row 5 (0-based data row 1)
Set-CellText $t 5 1 "84÷6="   # was "39÷9="
Set-CellText $t 5 2 "91÷2="   # was "68÷7="
Set-CellText $t 5 3 "14÷3="   # was "52÷9="
Set-CellText $t 5 4 "89÷2="   # was "31÷5="
Set-CellText $t 5 5 "44÷9="   # was "99÷4="

# Row 9 (0-based data row 2)
Set-CellText $t 9 1 "71÷4="   # was "11÷8="
Set-CellText $t 9 2 "10÷3="   # was "14÷4="
Set-CellText $t 9 3 "15÷6="   # was "57÷6="
Set-CellText $t 9 4 "20÷8="   # was "56÷8="
Set-CellText $t 9 5 "69÷2="   # was "56÷8="

# Row 13 (0-based data row 3)
Set-CellText $t 13 1 "40÷9="   # was "59÷2="
Set-CellText $t 13 2 "40÷8="   # was "25÷5="
Set-CellText $t 13 3 "19÷2="   # was "27÷7="
Set-CellText $t 13 4 "78÷9="   # was "41÷8="
Set-CellText $t 13 5 "82÷8="   # was "15÷5="

# Row 17 (0-based data row 4)
Set-CellText $t 17 1 "54÷9="   # was "81÷7="
Set-CellText $t 17 2 "10÷6="   # was "69÷3="
Set-CellText $t 17 3 "24÷5="   # was "53÷9="
Set-CellText $t 17 4 "17÷8="   # was "71÷7="
Set-CellText $t 17 5 "87÷9="   # was "58÷8="
